$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 146.620486
$ws.Range("H2").Value = 439.861458
$ws.Range("I2").Value = 0.3983053592962091
$ws.Range("J2").Value = 0.3983053592962091
$ws.Range("M2").Value = 16.23967033333334
$ws.Range("N2").Value = 48.71901100000001
$ws.Range("O2").Value = 0.3412424148893533
$ws.Range("P2").Value = 0.3412424148893533
$ws.Range("Q2").Value = 2381.068356753116
$ws.Range("R2").Value = 21429.61521077804
$ws.Range("S2").Value = 0.1359186826696099
$ws.Range("T2").Value = 0.1359186826696099
$ws.Range("G3").Value = 146.620486
$ws.Range("H3").Value = 439.861458
$ws.Range("I3").Value = 0.3983053592962091
$ws.Range("J3").Value = 0.3983053592962091
$ws.Range("O3").Value = 0.4874916916781935
$ws.Range("P3").Value = 0.4874916916781935
$ws.Range("Q3").Value = 3401.543860282909
$ws.Range("R3").Value = 30613.89474254618
$ws.Range("S3").Value = 0.1941705534077997
$ws.Range("T3").Value = 0.1941705534077996
$ws.Range("G4").Value = 146.620486
$ws.Range("H4").Value = 439.861458
$ws.Range("I4").Value = 0.3983053592962091
$ws.Range("J4").Value = 0.3983053592962091
$ws.Range("M4").Value = 8.150515666666667
$ws.Range("N4").Value = 24.451547
$ws.Range("O4").Value = 0.1712658934324533
$ws.Range("P4").Value = 0.1712658934324533
$ws.Range("Q4").Value = 1195.032568197281
$ws.Range("R4").Value = 10755.29311377553
$ws.Range("S4").Value = 0.06821612321879958
$ws.Range("T4").Value = 0.06821612321879957
$ws.Range("I5").Value = 0.534552907532962
$ws.Range("J5").Value = 0.5345529075329621
$ws.Range("M5").Value = 16.23967033333334
$ws.Range("N5").Value = 48.71901100000001
$ws.Range("O5").Value = 0.3412424148893533
$ws.Range("P5").Value = 0.3412424148893533
$ws.Range("Q5").Value = 3195.555830295915
$ws.Range("R5").Value = 28760.00247266324
$ws.Range("S5").Value = 0.1824121250526732
$ws.Range("T5").Value = 0.1824121250526732
$ws.Range("I6").Value = 0.534552907532962
$ws.Range("J6").Value = 0.5345529075329621
$ws.Range("O6").Value = 0.4874916916781935
$ws.Range("P6").Value = 0.4874916916781935
$ws.Range("S6").Value = 0.2605901011847406
$ws.Range("T6").Value = 0.2605901011847406
$ws.Range("I7").Value = 0.534552907532962
$ws.Range("J7").Value = 0.5345529075329621
$ws.Range("M7").Value = 8.150515666666667
$ws.Range("N7").Value = 24.451547
$ws.Range("O7").Value = 0.1712658934324533
$ws.Range("P7").Value = 0.1712658934324533
$ws.Range("Q7").Value = 1603.815060523388
$ws.Range("R7").Value = 14434.33554471049
$ws.Range("S7").Value = 0.09155068129554834
$ws.Range("T7").Value = 0.09155068129554836
$ws.Range("G8").Value = 24.174389
$ws.Range("H8").Value = 72.523167
$ws.Range("I8").Value = 0.0656715098899026
$ws.Range("J8").Value = 0.0656715098899026
$ws.Range("M8").Value = 16.23967033333334
$ws.Range("N8").Value = 48.71901100000001
$ws.Range("O8").Value = 0.3412424148893533
$ws.Range("P8").Value = 0.3412424148893533
$ws.Range("Q8").Value = 392.5841078697597
$ws.Range("R8").Value = 3533.256970827838
$ws.Range("S8").Value = 0.02240990462426041
$ws.Range("T8").Value = 0.02240990462426041
$ws.Range("G9").Value = 24.174389
$ws.Range("H9").Value = 72.523167
$ws.Range("I9").Value = 0.0656715098899026
$ws.Range("J9").Value = 0.0656715098899026
$ws.Range("O9").Value = 0.4874916916781935
$ws.Range("P9").Value = 0.4874916916781935
$ws.Range("Q9").Value = 560.8373476475906
$ws.Range("R9").Value = 5047.536128828316
$ws.Range("S9").Value = 0.03201431545128983
$ws.Range("T9").Value = 0.03201431545128983
$ws.Range("G10").Value = 24.174389
$ws.Range("H10").Value = 72.523167
$ws.Range("I10").Value = 0.0656715098899026
$ws.Range("J10").Value = 0.0656715098899026
$ws.Range("M10").Value = 8.150515666666667
$ws.Range("N10").Value = 24.451547
$ws.Range("O10").Value = 0.1712658934324533
$ws.Range("P10").Value = 0.1712658934324533
$ws.Range("Q10").Value = 197.0337362765944
$ws.Range("R10").Value = 1773.303626489349
$ws.Range("S10").Value = 0.01124728981435236
$ws.Range("T10").Value = 0.01124728981435236
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.541205
$ws.Range("H11").Value = 1.623615
$ws.Range("I11").Value = 0.001470223280926138
$ws.Range("J11").Value = 0.001470223280926138
$ws.Range("M11").Value = 16.23967033333334
$ws.Range("N11").Value = 48.71901100000001
$ws.Range("O11").Value = 0.3412424148893533
$ws.Range("P11").Value = 0.3412424148893533
$ws.Range("Q11").Value = 8.78899078275167
$ws.Range("R11").Value = 79.10091704476501
$ws.Range("S11").Value = 0.0005017025428097834
$ws.Range("T11").Value = 0.0005017025428097834
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.541205
$ws.Range("H12").Value = 1.623615
$ws.Range("I12").Value = 0.001470223280926138
$ws.Range("J12").Value = 0.001470223280926138
$ws.Range("O12").Value = 0.4874916916781935
$ws.Range("P12").Value = 0.4874916916781935
$ws.Range("Q12").Value = 12.55576621744667
$ws.Range("R12").Value = 113.00189595702
$ws.Range("S12").Value = 0.0007167216343633469
$ws.Range("T12").Value = 0.0007167216343633468
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.541205
$ws.Range("H13").Value = 1.623615
$ws.Range("I13").Value = 0.001470223280926138
$ws.Range("J13").Value = 0.001470223280926138
$ws.Range("M13").Value = 8.150515666666667
$ws.Range("N13").Value = 24.451547
$ws.Range("O13").Value = 0.1712658934324533
$ws.Range("P13").Value = 0.1712658934324533
$ws.Range("Q13").Value = 4.411099831378334
$ws.Range("R13").Value = 39.699898482405
$ws.Range("S13").Value = 0.0002517991037530078
$ws.Range("T13").Value = 0.0002517991037530078
